# Scheduled runner update: refresh market-price derived columns (H:N)
# across the per-profession Leve-profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1000.7143
$ws.Range("I70").Value = 935.8889
$ws.Range("J70").Value = 1219.5
$ws.Range("K70").Value = 2807.6667
$ws.Range("L70").Value = 3658.5
$ws.Range("M70").Value = -2537.6667
$ws.Range("N70").Value = -4198.5
$ws.Range("H73").Value = 1000.7143
$ws.Range("I73").Value = 935.8889
$ws.Range("J73").Value = 1219.5
$ws.Range("K73").Value = 2807.6667
$ws.Range("L73").Value = 3658.5
$ws.Range("M73").Value = -1871.6667
$ws.Range("N73").Value = -5530.5
$ws.Range("H76").Value = 2992.7856
$ws.Range("I76").Value = 2883.25
$ws.Range("J76").Value = 3650
$ws.Range("K76").Value = 2883.25
$ws.Range("L76").Value = 3650
$ws.Range("M76").Value = -2568.25
$ws.Range("N76").Value = -4280
$ws.Range("H79").Value = 2992.7856
$ws.Range("I79").Value = 2883.25
$ws.Range("J79").Value = 3650
$ws.Range("K79").Value = 2883.25
$ws.Range("L79").Value = 3650
$ws.Range("M79").Value = -1791.25
$ws.Range("N79").Value = -5834
$ws.Range("H86").Value = 202600.33
$ws.Range("I86").Value = 303000.75
$ws.Range("J86").Value = 1799.5
$ws.Range("K86").Value = 303000.75
$ws.Range("L86").Value = 1799.5
$ws.Range("M86").Value = -301877.75
$ws.Range("N86").Value = -4045.5
$ws.Range("H89").Value = 202600.33
$ws.Range("I89").Value = 303000.75
$ws.Range("J89").Value = 1799.5
$ws.Range("K89").Value = 1515003.75
$ws.Range("L89").Value = 8997.5
$ws.Range("M89").Value = -1509387.75
$ws.Range("N89").Value = -20229.5
$ws.Range("H113").Value = 3509.087
$ws.Range("I113").Value = 2618.1428
$ws.Range("J113").Value = 4895
$ws.Range("K113").Value = 2618.1428
$ws.Range("L113").Value = 4895
$ws.Range("M113").Value = 635.8571999999999
$ws.Range("N113").Value = -11403
$ws.Range("H132").Value = 5392.1113
$ws.Range("I132").Value = 4023.1904
$ws.Range("J132").Value = 10183.333
$ws.Range("K132").Value = 12069.5712
$ws.Range("L132").Value = 30549.999
$ws.Range("M132").Value = -9539.5712
$ws.Range("N132").Value = -35609.999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9764.134
$ws.Range("I32").Value = 8797.236000000001
$ws.Range("K32").Value = 8797.236000000001
$ws.Range("M32").Value = -8510.236000000001
$ws.Range("H61").Value = 2730.96
$ws.Range("I61").Value = 3312.4285
$ws.Range("J61").Value = 1990.909
$ws.Range("K61").Value = 3312.4285
$ws.Range("L61").Value = 1990.909
$ws.Range("M61").Value = -3100.4285
$ws.Range("N61").Value = -2414.909
$ws.Range("H88").Value = 1417.75
$ws.Range("I88").Value = 1356.2222
$ws.Range("J88").Value = 1602.3334
$ws.Range("K88").Value = 1356.2222
$ws.Range("L88").Value = 1602.3334
$ws.Range("M88").Value = -950.2221999999999
$ws.Range("N88").Value = -2414.3334
$ws.Range("H91").Value = 1417.75
$ws.Range("I91").Value = 1356.2222
$ws.Range("J91").Value = 1602.3334
$ws.Range("K91").Value = 1356.2222
$ws.Range("L91").Value = 1602.3334
$ws.Range("M91").Value = 47.77780000000007
$ws.Range("N91").Value = -4410.3334
$ws.Range("H97").Value = 778.5599999999999
$ws.Range("I97").Value = 778.5599999999999
$ws.Range("K97").Value = 778.5599999999999
$ws.Range("M97").Value = -282.5599999999999
$ws.Range("H132").Value = 703823.9
$ws.Range("I132").Value = 1702778.6
$ws.Range("J132").Value = 4555.5
$ws.Range("K132").Value = 5108335.800000001
$ws.Range("L132").Value = 13666.5
$ws.Range("M132").Value = -5105805.800000001
$ws.Range("N132").Value = -18726.5
$ws.Range("H136").Value = 2730.96
$ws.Range("I136").Value = 3312.4285
$ws.Range("J136").Value = 1990.909
$ws.Range("K136").Value = 9937.2855
$ws.Range("L136").Value = 5972.727000000001
$ws.Range("M136").Value = -7387.2855
$ws.Range("N136").Value = -11072.727

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1364.5625
$ws.Range("I86").Value = 1080.5
$ws.Range("J86").Value = 1535
$ws.Range("K86").Value = 1080.5
$ws.Range("L86").Value = 1535
$ws.Range("M86").Value = 42.5
$ws.Range("N86").Value = -3781
$ws.Range("H89").Value = 1364.5625
$ws.Range("I89").Value = 1080.5
$ws.Range("J89").Value = 1535
$ws.Range("K89").Value = 5402.5
$ws.Range("L89").Value = 7675
$ws.Range("M89").Value = 213.5
$ws.Range("N89").Value = -18907
$ws.Range("H94").Value = 938.6539
$ws.Range("I94").Value = 848.25
$ws.Range("K94").Value = 848.25
$ws.Range("M94").Value = -397.25
$ws.Range("H105").Value = 1865.8
$ws.Range("I105").Value = 1485
$ws.Range("J105").Value = 2437
$ws.Range("K105").Value = 1485
$ws.Range("L105").Value = 2437
$ws.Range("M105").Value = 262
$ws.Range("N105").Value = -5931
$ws.Range("H134").Value = 4441.418
$ws.Range("I134").Value = 2165.2
$ws.Range("J134").Value = 7172.88
$ws.Range("K134").Value = 6495.599999999999
$ws.Range("L134").Value = 21518.64
$ws.Range("M134").Value = -3960.599999999999
$ws.Range("N134").Value = -26588.64

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2754.4375
$ws.Range("I31").Value = 2448.4614
$ws.Range("J31").Value = 2963.7896
$ws.Range("K31").Value = 2448.4614
$ws.Range("L31").Value = 2963.7896
$ws.Range("M31").Value = -2153.4614
$ws.Range("N31").Value = -3553.7896
$ws.Range("H34").Value = 2754.4375
$ws.Range("I34").Value = 2448.4614
$ws.Range("J34").Value = 2963.7896
$ws.Range("K34").Value = 2448.4614
$ws.Range("L34").Value = 2963.7896
$ws.Range("M34").Value = -2246.4614
$ws.Range("N34").Value = -3367.7896
$ws.Range("H132").Value = 2403.7942
$ws.Range("I132").Value = 1583.4348
$ws.Range("J132").Value = 4119.091
$ws.Range("K132").Value = 4750.3044
$ws.Range("L132").Value = 12357.273
$ws.Range("M132").Value = -2220.3044
$ws.Range("N132").Value = -17417.273
$ws.Range("H134").Value = 1869.0834
$ws.Range("I134").Value = 1386.9445
$ws.Range("J134").Value = 3315.5
$ws.Range("K134").Value = 4160.833500000001
$ws.Range("L134").Value = 9946.5
$ws.Range("M134").Value = -1625.833500000001
$ws.Range("N134").Value = -15016.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 505
$ws.Range("I57").Value = 505
$ws.Range("K57").Value = 1515
$ws.Range("M57").Value = -956

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5982.5713
$ws.Range("I70").Value = 5363.143
$ws.Range("J70").Value = 6602
$ws.Range("K70").Value = 5363.143
$ws.Range("L70").Value = 6602
$ws.Range("M70").Value = -5093.143
$ws.Range("N70").Value = -7142
$ws.Range("H73").Value = 5982.5713
$ws.Range("I73").Value = 5363.143
$ws.Range("J73").Value = 6602
$ws.Range("K73").Value = 5363.143
$ws.Range("L73").Value = 6602
$ws.Range("M73").Value = -4427.143
$ws.Range("N73").Value = -8474
$ws.Range("H113").Value = 1912.5
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 270
$ws.Range("N113").Value = -6340

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 38240.242
$ws.Range("I132").Value = 53398.5
$ws.Range("J132").Value = 4555.222
$ws.Range("K132").Value = 160195.5
$ws.Range("L132").Value = 13665.666
$ws.Range("M132").Value = -157665.5
$ws.Range("N132").Value = -18725.666
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6450
$ws.Range("N136").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1780.84
$ws.Range("I132").Value = 1379.7097
$ws.Range("J132").Value = 2435.3157
$ws.Range("K132").Value = 4139.1291
$ws.Range("L132").Value = 7305.9471
$ws.Range("M132").Value = -1609.1291
$ws.Range("N132").Value = -12365.9471
$ws.Range("H136").Value = 1374072.1
$ws.Range("I136").Value = 3706024.2
$ws.Range("J136").Value = 2335.5881
$ws.Range("K136").Value = 11118072.6
$ws.Range("L136").Value = 7006.7643
$ws.Range("M136").Value = -11115522.6
$ws.Range("N136").Value = -12106.7643
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
